$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new product row (row 36) below the existing data (1:35)
$ws.Range("A36").Value = 2943549
$ws.Range("B36").Value = "Pril Power blue 650 ml"
$ws.Range("C36").Value = 12
$ws.Range("D36").Value = 168

# Match the author's on-screen view: scrolled a couple of rows back up
# and the active selection moved to F30.
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F30").Select()
